$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove "20191411001" from A2 (keep its quote-prefixed formatting),
# add "Vivek" to B2
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "Vivek"

# Row 3: A3 is fully cleared (content + formatting); B3/C3 stay as-is
# (Vivek / Singla); D3 changes from "v" to "*"
$ws.Range("A3").Clear()
$ws.Range("D3").Value = "*"

# Row 4: B4/C4/D4 ("Vivek"/"Singla"/"*") are cleared; A4 becomes "*"
$ws.Range("B4").Clear()
$ws.Range("C4").Clear()
$ws.Range("D4").Clear()
$ws.Range("A4").Value = "*"

# Row 5: A5 ("*") is cleared; B5 and C5 become "*"
$ws.Range("A5").Clear()
$ws.Range("B5").Value = "*"
$ws.Range("C5").Value = "*"

# Row 6: B6/C6 ("*") are cleared; A6 becomes "20191411001" (quote-prefixed,
# matching the original formatting used for that value)
$ws.Range("B6").Clear()
$ws.Range("C6").Clear()
$ws.Range("A6").Value = "'20191411001"

# Update the active selection to C10
$ws.Range("C10").Select()
